$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '48.048.21'
Set-TextValue $ws.Range('E2') '  +0.48%  '
Set-TextValue $ws.Range('D3') '2.500.48'
Set-TextValue $ws.Range('E3') '  +0.06%  '
Set-TextValue $ws.Range('D4') '0.998'
Set-TextValue $ws.Range('E4') '  -0.24%  '
Set-TextValue $ws.Range('D5') '320.36'
Set-TextValue $ws.Range('E5') '  -0.82%  '
Set-TextValue $ws.Range('D6') '107.36'
Set-TextValue $ws.Range('E6') '  -1.30%  '
Set-TextValue $ws.Range('D7') '0.525'
Set-TextValue $ws.Range('E7') '  +0.20%  '
Set-TextValue $ws.Range('D8') '0.998'
Set-TextValue $ws.Range('E8') '  -0.17%  '
Set-TextValue $ws.Range('D9') '0.541'
Set-TextValue $ws.Range('E9') '  -1.72%  '
Set-TextValue $ws.Range('D10') '39.66'
Set-TextValue $ws.Range('E10') '  -1.31%  '
Set-TextValue $ws.Range('D11') '20.13'
Set-TextValue $ws.Range('E11') '  +6.58%  '
Set-TextValue $ws.Range('E12') '  -0.37%  '
Set-TextValue $ws.Range('E13') '  -0.11%  '
Set-TextValue $ws.Range('D14') '7.10'
Set-TextValue $ws.Range('E14') '  -1.60%  '
Set-TextValue $ws.Range('D15') '2.890.95'
Set-TextValue $ws.Range('E15') '  -0.05%  '
Set-TextValue $ws.Range('D16') '2.501.88'
Set-TextValue $ws.Range('E16') '  +0.38%  '
Set-TextValue $ws.Range('D17') '0.834'
Set-TextValue $ws.Range('E17') '  -1.98%  '
Set-TextValue $ws.Range('D18') '47.880.29'
Set-TextValue $ws.Range('E18') '  +0.34%  '
Set-TextValue $ws.Range('D19') '12.93'
Set-TextValue $ws.Range('E19') '  -1.62%  '
Set-TextValue $ws.Range('D20') '6.69'
Set-TextValue $ws.Range('E20') '  +0.99%  '
Set-TextValue $ws.Range('D21') '0.0₃0940'
Set-TextValue $ws.Range('E21') '  -0.07%  '
Set-TextValue $ws.Range('D22') '2.76'
Set-TextValue $ws.Range('E22') '  -0.44%  '
Set-TextValue $ws.Range('D23') '275.76'
Set-TextValue $ws.Range('E23') '  +11.20%  '
Set-TextValue $ws.Range('D24') '71.48'
Set-TextValue $ws.Range('E24') '  +1.00%  '
Set-TextValue $ws.Range('E25') '  -1.31%  '
Set-TextValue $ws.Range('E26') '  -0.08%  '
Set-TextValue $ws.Range('D27') '25.85'
Set-TextValue $ws.Range('E27') '  -0.24%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D28') '2.26'
Set-TextValue $ws.Range('E28') '  +2.79%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D29') '9.70'
Set-TextValue $ws.Range('E29') '  -2.70%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D30') '0.141'
Set-TextValue $ws.Range('E30') '  +1.54%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D31') '35.13'
Set-TextValue $ws.Range('E31') '  -0.12%  '
Set-TextValue $ws.Range('D32') '49.71'
Set-TextValue $ws.Range('E32') '  -0.34%  '
Set-TextValue $ws.Range('D33') '19.45'
Set-TextValue $ws.Range('E33') '  -2.06%  '
Set-TextValue $ws.Range('E34') '  -0.21%  '
Set-TextValue $ws.Range('E35') '  -1.06%  '
Set-TextValue $ws.Range('E36') '  -0.65%  '
Set-TextValue $ws.Range('D37') '1.94'
Set-TextValue $ws.Range('E37') '  -1.03%  '
Set-TextValue $ws.Range('D38') '4.62'
Set-TextValue $ws.Range('E38') '  -1.04%  '
Set-TextValue $ws.Range('E39') '  -3.14%  '
Set-TextValue $ws.Range('E40') '  -0.61%  '
Set-TextValue $ws.Range('D41') '121.36'
Set-TextValue $ws.Range('E41') '  +1.81%  '
Set-TextValue $ws.Range('E42') '  -0.28%  '
Set-TextValue $ws.Range('D43') '21.29'
Set-TextValue $ws.Range('E43') '  -3.85%  '
Set-TextValue $ws.Range('D44') '0.0303'
Set-TextValue $ws.Range('E44') '  +1.70%  '
Set-TextValue $ws.Range('D45') '2.017.34'
Set-TextValue $ws.Range('E45') '  +0.85%  '
Set-TextValue $ws.Range('D46') '3.13'
Set-TextValue $ws.Range('E46') '  +2.60%  '
Set-TextValue $ws.Range('D47') '2.00'
Set-TextValue $ws.Range('E47') '  -1.56%  '
Set-TextValue $ws.Range('D48') '1.84'
Set-TextValue $ws.Range('E48') '  +1.75%  '
Set-TextValue $ws.Range('D49') '8.99'
Set-TextValue $ws.Range('E49') '  -0.42%  '
Set-TextValue $ws.Range('E50') '  +1.01%  '
Set-TextValue $ws.Range('D51') '80.40'
Set-TextValue $ws.Range('E51') '  +3.36%  '
